$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Update the query text held in the B/C columns (rows 2-4) with the new
#    Cypher queries described in the commit.
# ---------------------------------------------------------------------------

$statsQuery = @"
CALL{
    MATCH (p:participant)-->(s:study)
    OPTIONAL MATCH (samp:sample)-->(p)
    OPTIONAL MATCH (samp)<--(f:file)
    OPTIONAL MATCH (f)<--(g:genomic_info)
    OPTIONAL MATCH (p)<--(diag:diagnosis)
    WITH s, p, samp, f, g, diag, apoc.coll.flatten(COLLECT (apoc.text.split(f.experimental_strategy_and_data_subtypes,"[;,]\s{0,1}")), true) as es
    WHERE "Methylation Array" IN es
    RETURN 
        count(distinct p) AS num_participants
}
WITH num_participants
CALL {
    MATCH (samp:sample)-->(p:participant)-->(s:study)
    OPTIONAL MATCH (samp)<--(f:file)
    OPTIONAL MATCH (p)<--(diag:diagnosis)
    OPTIONAL MATCH (f)<--(g:genomic_info)
    OPTIONAL MATCH (p)<--(diag:diagnosis)
    WITH s, p, samp, f, g, diag, apoc.coll.flatten(COLLECT (apoc.text.split(f.experimental_strategy_and_data_subtypes,"[;,]\s{0,1}")), true) as es
    WHERE "Methylation Array" IN es
    RETURN 
        count(distinct samp) AS num_samples
}
WITH num_participants, num_samples
CALL {
    MATCH (f:file)-->(s:study)
    OPTIONAL MATCH (f)<--(g:genomic_info)
    OPTIONAL MATCH (samp:sample)<--(f)
    OPTIONAL MATCH (p:participant)<--(samp)
    OPTIONAL MATCH (p)<--(diag:diagnosis)
    WITH s, p, samp, f, g, diag, apoc.coll.flatten(COLLECT (apoc.text.split(f.experimental_strategy_and_data_subtypes,"[;,]\s{0,1}")), true) as es
    WHERE "Methylation Array" IN es
    RETURN 
        count(distinct s) AS num_studies,
        count(distinct f) AS num_files
}
RETURN 
    num_studies AS Studies,
    num_participants AS Participants,
    num_samples AS Samples,
    num_files AS ``Files``
"@

$participantsQuery = @"
MATCH (p:participant)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
OPTIONAL MATCH (samp)<--(f:file)
WITH p, samp, apoc.coll.flatten(COLLECT (apoc.text.split(f.experimental_strategy_and_data_subtypes,"[;,]\s{0,1}")), true) as es
WHERE "Methylation Array" IN es
WITH p
OPTIONAL MATCH (p)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
WITH s, p, apoc.coll.sort(collect(distinct coalesce(samp.sample_id, "Not specified in data"))) as samp
RETURN 
coalesce(p.participant_id,'') as ``Participant ID``,
coalesce(s.study_name, '') as ``Study Name``,
coalesce(s.phs_accession,'') as ``Accession``,
coalesce(p.gender,'') as ``Gender``,
coalesce(apoc.text.join(samp, ','), '') as ``Samples``
ORDER BY p.participant_id
LIMIT 100
"@

$samplesQuery = @"
MATCH (samp:sample)-->(p:participant)-->(s:study)
OPTIONAL MATCH (samp)<--(f:file)
OPTIONAL MATCH (f)<--(g:genomic_info)
OPTIONAL MATCH (p)<--(diag:diagnosis)
WITH s, p, samp, f, g, diag, apoc.coll.flatten(COLLECT (apoc.text.split(f.experimental_strategy_and_data_subtypes,"[;,]\s{0,1}")), true) as es
WHERE "Methylation Array" IN es
WITH DISTINCT s, p, samp
RETURN
    coalesce(samp.sample_id, '') as ``Sample ID``,
    coalesce(p.participant_id,'') as ``Participant ID``,
    coalesce(s.study_name, '') as ``Study Name``,
    coalesce(s.phs_accession,'') as ``Accession``,
    coalesce(samp.sample_tumor_status,'') as ``Tumor``,
    coalesce(samp.sample_type,'') as ``Analyte Type``
ORDER BY samp.sample_id LIMIT 100
"@

$filesQuery = @"
MATCH (f:file)-->(s:study)
OPTIONAL MATCH (samp:sample)<--(f)
OPTIONAL MATCH (samp)-->(p:participant)
OPTIONAL MATCH (f)<--(g:genomic_info)
OPTIONAL MATCH (p)<--(diag:diagnosis)
WITH s, p, samp, f, g, diag, apoc.coll.flatten(COLLECT (apoc.text.split(f.experimental_strategy_and_data_subtypes,"[;,]\s{0,1}")), true) as es
WHERE "Methylation Array" IN es
WITH DISTINCT f, s, p, samp
RETURN
    coalesce(f.file_name, '') as ``File Name``,
    coalesce(s.study_name,'') as ``Study Name``,
    coalesce(s.phs_accession,'') as ``Accession``,
    coalesce(p.participant_id, '') as ``Participant ID``,
    coalesce(samp.sample_id, '') as ``Sample ID``,
    coalesce(f.file_type, '') as ``File Type``
ORDER BY f.file_name LIMIT 100
"@

$ws.Range("C2").Value = $statsQuery
$ws.Range("B2").Value = $participantsQuery
$ws.Range("C3").Value = $statsQuery
$ws.Range("B3").Value = $samplesQuery
$ws.Range("C4").Value = $statsQuery
$ws.Range("B4").Value = $filesQuery

# ---------------------------------------------------------------------------
# 2. Bump the whole sheet's font size from 11/12pt to 14pt.
# ---------------------------------------------------------------------------

$ws.Cells.Font.Size = 14

# ---------------------------------------------------------------------------
# 3. Update the view: scroll so row 3 is at the top, and select E4.
# ---------------------------------------------------------------------------

$ws.Range("E4").Select()
$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 1
